$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old country names -> new country names (column A, "COUNTRY")
$renames = @{
    "Bahamas, The"                        = "Bahamas"
    "Congo, Democratic Republic of the"   = "Congo (Kinshasa)"
    "Congo, Republic of the"              = "Congo (Brazzaville)"
    "Czech Republic"                      = "Czechia"
    "Gambia, The"                         = "Gambia"
    "Swaziland"                           = "Eswatini"
    "Taiwan"                              = "Taiwan*"
    "United States"                       = "US"
    "West Bank"                           = "West Bank and Gaza"
}

$colA = $ws.Range("A1:A" + $ws.UsedRange.Rows.Count)

foreach ($oldName in $renames.Keys) {
    $newName = $renames[$oldName]
    $cell = $colA.Find($oldName, [Type]::Missing, [Type]::Missing, 1)
    if ($cell -ne $null) {
        $cell.Value = $newName
    }
}
